$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the cells that must disappear entirely (content + formatting) ---
$ws.Range("D1:D4").Clear()
$ws.Range("C2:C4").Clear()
$ws.Range("E6").Clear()

# --- Row 1 (headers) ---
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "email"
$ws.Range("C1").Value = "datos"

# --- Row 2 ---
$ws.Range("A2").Value = "Julian 1"
$ws.Range("B2").Value = "sd"

# --- Row 3 ---
$ws.Range("A3").Value = "Julian 2"
$ws.Range("B3").Value = "asd"

# --- Row 4 ---
$ws.Range("A4").Value = "Julian 4"
$ws.Range("B4").Value = "aaadfs"

# --- Row 5 (new row) ---
$ws.Range("A5").Value = "julian 4"
$ws.Range("B5").Value = "sss"
$ws.Range("C5").Value = "aaa"

# --- Selection matches the saved cursor position from the diff ---
$ws.Range("C5").Select()
